$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 ("dimension/measure" row),
# pushing the existing rows 2-5 down to rows 3-6. The new row 2 holds the
# machine-friendly ("slug") column keys that let two columns be related
# to build a SKOS hierarchy.
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "ccaa-nombre"
$ws.Range("B2").Value = "renta-disponible-bruta"
$ws.Range("C2").Value = "codigo"
$ws.Range("D2").Value = "comarca"
$ws.Range("E2").Value = "ccaa-codigo"
$ws.Range("F2").Value = "municipio-codigo"
$ws.Range("G2").Value = "ano"
$ws.Range("H2").Value = "renta-disponible-bruta-per-capita"
$ws.Range("I2").Value = "municipio-nombre"

# (The newly inserted row already inherits the "s=1" cell format used by
# the rest of the sheet, matching the row above it.)

# The trailing stray "mapping-ano.xlsx" row (old row 5, now shifted to
# row 6 by the insert above) is obsolete cruft and is dropped entirely.
$ws.Rows.Item(6).Delete()
